$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$data[0,0] = "ECs"
$data[0,1] = "Il6"
$data[0,2] = "Il6st"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 7.292895000000001
$data[0,7] = 21.878685
$data[0,8] = 0.2565758520803378
$data[0,9] = 0.2565758520803378
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 39.96608766666667
$data[0,13] = 119.898263
$data[0,14] = 0.2616165719423124
$data[0,15] = 0.2616165719423124
$data[0,16] = 291.4684809137951
$data[0,17] = 2623.216328224155
$data[0,18] = 0.06712449486443579
$data[0,19] = 0.06712449486443579
$data[1,0] = "ECs"
$data[1,1] = "Il6"
$data[1,2] = "Il6st"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 7.292895000000001
$data[1,7] = 21.878685
$data[1,8] = 0.2565758520803378
$data[1,9] = 0.2565758520803378
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 85.11961100000001
$data[1,13] = 255.358833
$data[1,14] = 0.5571899111219771
$data[1,15] = 0.557189911121977
$data[1,16] = 620.7683854638451
$data[1,17] = 5586.915469174605
$data[1,18] = 0.1429614762166889
$data[1,19] = 0.1429614762166889
$data[2,0] = "ECs"
$data[2,1] = "Il6"
$data[2,2] = "Il6st"
$data[2,3] = "sCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 7.292895000000001
$data[2,7] = 21.878685
$data[2,8] = 0.2565758520803378
$data[2,9] = 0.2565758520803378
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 27.68018833333333
$data[2,13] = 83.040565
$data[2,14] = 0.1811935169357105
$data[2,15] = 0.1811935169357105
$data[2,16] = 201.868707095225
$data[2,17] = 1816.818363857025
$data[2,18] = 0.04648988099921302
$data[2,19] = 0.04648988099921302
$data[3,0] = "FAPs"
$data[3,1] = "Il6"
$data[3,2] = "Il6st"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 14.46476266666667
$data[3,7] = 43.394288
$data[3,8] = 0.5088937666509471
$data[3,9] = 0.5088937666509471
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 39.96608766666667
$data[3,13] = 119.898263
$data[3,14] = 0.2616165719423124
$data[3,15] = 0.2616165719423124
$data[3,16] = 578.0999728135272
$data[3,17] = 5202.899755321745
$data[3,18] = 0.1331350427140318
$data[3,19] = 0.1331350427140318
$data[4,0] = "FAPs"
$data[4,1] = "Il6"
$data[4,2] = "Il6st"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 14.46476266666667
$data[4,7] = 43.394288
$data[4,8] = 0.5088937666509471
$data[4,9] = 0.5088937666509471
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 85.11961100000001
$data[4,13] = 255.358833
$data[4,14] = 0.5571899111219771
$data[4,15] = 0.557189911121977
$data[4,16] = 1231.23497139399
$data[4,17] = 11081.11474254591
$data[4,18] = 0.2835504726107694
$data[4,19] = 0.2835504726107693
$data[5,0] = "FAPs"
$data[5,1] = "Il6"
$data[5,2] = "Il6st"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 14.46476266666667
$data[5,7] = 43.394288
$data[5,8] = 0.5088937666509471
$data[5,9] = 0.5088937666509471
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 27.68018833333333
$data[5,13] = 83.040565
$data[5,14] = 0.1811935169357105
$data[5,15] = 0.1811935169357105
$data[5,16] = 400.3873548103022
$data[5,17] = 3603.48619329272
$data[5,18] = 0.09220825132614587
$data[5,19] = 0.09220825132614587
$data[6,0] = "sCs"
$data[6,1] = "Il6"
$data[6,2] = "Il6st"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 6.666276
$data[6,7] = 19.998828
$data[6,8] = 0.2345303812687151
$data[6,9] = 0.2345303812687151
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 39.96608766666667
$data[6,13] = 119.898263
$data[6,14] = 0.2616165719423124
$data[6,15] = 0.2616165719423124
$data[6,16] = 266.424971026196
$data[6,17] = 2397.824739235764
$data[6,18] = 0.06135703436384476
$data[6,19] = 0.06135703436384475
$data[7,0] = "sCs"
$data[7,1] = "Il6"
$data[7,2] = "Il6st"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 6.666276
$data[7,7] = 19.998828
$data[7,8] = 0.2345303812687151
$data[7,9] = 0.2345303812687151
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 85.11961100000001
$data[7,13] = 255.358833
$data[7,14] = 0.5571899111219771
$data[7,15] = 0.557189911121977
$data[7,16] = 567.430819938636
$data[7,17] = 5106.877379447724
$data[7,18] = 0.1306779622945188
$data[7,19] = 0.1306779622945187
$data[8,0] = "sCs"
$data[8,1] = "Il6"
$data[8,2] = "Il6st"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 6.666276
$data[8,7] = 19.998828
$data[8,8] = 0.2345303812687151
$data[8,9] = 0.2345303812687151
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 27.68018833333333
$data[8,13] = 83.040565
$data[8,14] = 0.1811935169357105
$data[8,15] = 0.1811935169357105
$data[8,16] = 184.52377516198
$data[8,17] = 1660.71397645782
$data[8,18] = 0.04249538461035156
$data[8,19] = 0.04249538461035155

$ws.Range("A2:T10").Value2 = $data
